# Generate Report for Handoff
# Updates the "b.md" row (row 3) across all three sheets to reflect that the
# file is now ready for handoff with a new handoff package (commit hash
# 63290e5768f688058c7b37413b0a5c26c308f864) and updated timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: update Status (B3, C3) and Latest Handoff Date (D3) for b.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 06:35:45"

# ---------------------------------------------------------------------------
# zh-cn sheet: update Status (C3), Latest Handoff File (D3) and its
# hyperlink display text, and Latest Handoff Datetime (E3) for b.md
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-22 06:35:41"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet: update Status (C3), Latest Handoff File (D3) and its
# hyperlink display text, and Latest Handoff Datetime (E3) for b.md
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-22 06:35:45"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
